{"js": "// Auto-generated: replace each paragraph's text by position (document order),\n// matching the 101 <w:t> runs changed by the commit (the date line + 100\n// table-cell arithmetic expressions). Positional replacement is required\n// because some old values (e.g. \"69-9=\") occur more than once but map to\n// different new values depending on position.\nconst replacements = [\n  [\"2024-04-11 Thursday\", \"2024-04-12 Friday\"],\n  [\"29-26=\", \"97-23=\"],\n  [\"12-3=\", \"63-0=\"],\n  [\"29+33=\", \"20+11=\"],\n  [\"47-39=\", \"66-57=\"],\n  [\"97-89=\", \"26-21=\"],\n  [\"71-17=\", \"52+0=\"],\n  [\"88-61=\", \"94-82=\"],\n  [\"68-63=\", \"57-46=\"],\n  [\"48-31=\", \"11+55=\"],\n  [\"22+23=\", \"8+76=\"],\n  [\"31+37=\", \"63+3=\"],\n  [\"7+2=\", \"90-43=\"],\n  [\"51-2=\", \"64-10=\"],\n  [\"40+12=\", \"62-24=\"],\n  [\"12+47=\", \"12+78=\"],\n  [\"17+56=\", \"35+56=\"],\n  [\"83-81=\", \"64-41=\"],\n  [\"84-47=\", \"1+78=\"],\n  [\"96-11=\", \"87-44=\"],\n  [\"42+39=\", \"96-37=\"],\n  [\"18+80=\", \"45-21=\"],\n  [\"34-7=\", \"49+6=\"],\n  [\"67+18=\", \"74-13=\"],\n  [\"11+82=\", \"41-39=\"],\n  [\"52+23=\", \"3+70=\"],\n  [\"65-55=\", \"51-34=\"],\n  [\"19-4=\", \"39+60=\"],\n  [\"59-6=\", \"97-93=\"],\n  [\"83-54=\", \"59-51=\"],\n  [\"74-28=\", \"99-98=\"],\n  [\"85+13=\", \"34-24=\"],\n  [\"99-36=\", \"8+85=\"],\n  [\"22+56=\", \"91-35=\"],\n  [\"89+5=\", \"9-0=\"],\n  [\"34-5=\", \"51-45=\"],\n  [\"23+53=\", \"84-35=\"],\n  [\"42+12=\", \"21+56=\"],\n  [\"89-84=\", \"22-13=\"],\n  [\"44+18=\", \"66-41=\"],\n  [\"73-36=\", \"93-66=\"],\n  [\"3+65=\", \"23+68=\"],\n  [\"94-58=\", \"45+46=\"],\n  [\"37+27=\", \"14-1=\"],\n  [\"56-54=\", \"63+21=\"],\n  [\"80+3=\", \"77+20=\"],\n  [\"80+19=\", \"10+58=\"],\n  [\"71-15=\", \"22-4=\"],\n  [\"95-51=\", \"82-52=\"],\n  [\"8+86=\", \"48+22=\"],\n  [\"69-9=\", \"71-56=\"],\n  [\"36-36=\", \"88-14=\"],\n  [\"38-35=\", \"66-65=\"],\n  [\"93-53=\", \"91+1=\"],\n  [\"31+52=\", \"24+10=\"],\n  [\"66-45=\", \"12+46=\"],\n  [\"40+14=\", \"95-16=\"],\n  [\"75-11=\", \"60-39=\"],\n  [\"63-61=\", \"85-32=\"],\n  [\"16+34=\", \"82-55=\"],\n  [\"75+5=\", \"15+52=\"],\n  [\"46+21=\", \"7+17=\"],\n  [\"82-77=\", \"19+48=\"],\n  [\"54+20=\", \"48+6=\"],\n  [\"12+50=\", \"71+1=\"],\n  [\"49-42=\", \"80-20=\"],\n  [\"83+5=\", \"29+35=\"],\n  [\"34+56=\", \"9+37=\"],\n  [\"24+1=\", \"80-12=\"],\n  [\"39-11=\", \"91-55=\"],\n  [\"66-40=\", \"5+41=\"],\n  [\"74-30=\", \"78+0=\"],\n  [\"30+8=\", \"12+80=\"],\n  [\"25+69=\", \"80-62=\"],\n  [\"42+22=\", \"75+24=\"],\n  [\"9+68=\", \"93-45=\"],\n  [\"30+13=\", \"31+1=\"],\n  [\"47+28=\", \"28+71=\"],\n  [\"23+70=\", \"45-28=\"],\n  [\"52+43=\", \"40+41=\"],\n  [\"3+93=\", \"34+61=\"],\n  [\"69-9=\", \"41-25=\"],\n  [\"35+10=\", \"17+72=\"],\n  [\"18-16=\", \"48-45=\"],\n  [\"44-38=\", \"54+14=\"],\n  [\"11+88=\", \"7+62=\"],\n  [\"41+42=\", \"50+30=\"],\n  [\"69-21=\", \"19+65=\"],\n  [\"53+10=\", \"58+20=\"],\n  [\"63+1=\", \"53+22=\"],\n  [\"34-30=\", \"19+50=\"],\n  [\"17-11=\", \"40+44=\"],\n  [\"83-59=\", \"61-47=\"],\n  [\"59-35=\", \"89-53=\"],\n  [\"92-74=\", \"75-0=\"],\n  [\"44+22=\", \"72-8=\"],\n  [\"80-2=\", \"63+29=\"],\n  [\"17+63=\", \"11+68=\"],\n  [\"46-37=\", \"91+5=\"],\n  [\"3+19=\", \"48+31=\"],\n  [\"19+39=\", \"61-46=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  const current = para.text.trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Paragraph ${i} mismatch: expected \"${oldText}\", found \"${current}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated: replace each non-empty paragraph's text by position\n# (document order), matching the 101 text runs changed by the commit (the\n# date line + 100 table-cell arithmetic expressions). The Word COM\n# Paragraphs collection also yields one blank \"row end\" paragraph per\n# table row, so those are skipped. Positional replacement is required\n# because some old values (e.g. \"69-9=\") occur more than once but map to\n# different new values depending on position.\n$replacements = @(\n    @(\"2024-04-11 Thursday\", \"2024-04-12 Friday\"),\n    @(\"29-26=\", \"97-23=\"),\n    @(\"12-3=\", \"63-0=\"),\n    @(\"29+33=\", \"20+11=\"),\n    @(\"47-39=\", \"66-57=\"),\n    @(\"97-89=\", \"26-21=\"),\n    @(\"71-17=\", \"52+0=\"),\n    @(\"88-61=\", \"94-82=\"),\n    @(\"68-63=\", \"57-46=\"),\n    @(\"48-31=\", \"11+55=\"),\n    @(\"22+23=\", \"8+76=\"),\n    @(\"31+37=\", \"63+3=\"),\n    @(\"7+2=\", \"90-43=\"),\n    @(\"51-2=\", \"64-10=\"),\n    @(\"40+12=\", \"62-24=\"),\n    @(\"12+47=\", \"12+78=\"),\n    @(\"17+56=\", \"35+56=\"),\n    @(\"83-81=\", \"64-41=\"),\n    @(\"84-47=\", \"1+78=\"),\n    @(\"96-11=\", \"87-44=\"),\n    @(\"42+39=\", \"96-37=\"),\n    @(\"18+80=\", \"45-21=\"),\n    @(\"34-7=\", \"49+6=\"),\n    @(\"67+18=\", \"74-13=\"),\n    @(\"11+82=\", \"41-39=\"),\n    @(\"52+23=\", \"3+70=\"),\n    @(\"65-55=\", \"51-34=\"),\n    @(\"19-4=\", \"39+60=\"),\n    @(\"59-6=\", \"97-93=\"),\n    @(\"83-54=\", \"59-51=\"),\n    @(\"74-28=\", \"99-98=\"),\n    @(\"85+13=\", \"34-24=\"),\n    @(\"99-36=\", \"8+85=\"),\n    @(\"22+56=\", \"91-35=\"),\n    @(\"89+5=\", \"9-0=\"),\n    @(\"34-5=\", \"51-45=\"),\n    @(\"23+53=\", \"84-35=\"),\n    @(\"42+12=\", \"21+56=\"),\n    @(\"89-84=\", \"22-13=\"),\n    @(\"44+18=\", \"66-41=\"),\n    @(\"73-36=\", \"93-66=\"),\n    @(\"3+65=\", \"23+68=\"),\n    @(\"94-58=\", \"45+46=\"),\n    @(\"37+27=\", \"14-1=\"),\n    @(\"56-54=\", \"63+21=\"),\n    @(\"80+3=\", \"77+20=\"),\n    @(\"80+19=\", \"10+58=\"),\n    @(\"71-15=\", \"22-4=\"),\n    @(\"95-51=\", \"82-52=\"),\n    @(\"8+86=\", \"48+22=\"),\n    @(\"69-9=\", \"71-56=\"),\n    @(\"36-36=\", \"88-14=\"),\n    @(\"38-35=\", \"66-65=\"),\n    @(\"93-53=\", \"91+1=\"),\n    @(\"31+52=\", \"24+10=\"),\n    @(\"66-45=\", \"12+46=\"),\n    @(\"40+14=\", \"95-16=\"),\n    @(\"75-11=\", \"60-39=\"),\n    @(\"63-61=\", \"85-32=\"),\n    @(\"16+34=\", \"82-55=\"),\n    @(\"75+5=\", \"15+52=\"),\n    @(\"46+21=\", \"7+17=\"),\n    @(\"82-77=\", \"19+48=\"),\n    @(\"54+20=\", \"48+6=\"),\n    @(\"12+50=\", \"71+1=\"),\n    @(\"49-42=\", \"80-20=\"),\n    @(\"83+5=\", \"29+35=\"),\n    @(\"34+56=\", \"9+37=\"),\n    @(\"24+1=\", \"80-12=\"),\n    @(\"39-11=\", \"91-55=\"),\n    @(\"66-40=\", \"5+41=\"),\n    @(\"74-30=\", \"78+0=\"),\n    @(\"30+8=\", \"12+80=\"),\n    @(\"25+69=\", \"80-62=\"),\n    @(\"42+22=\", \"75+24=\"),\n    @(\"9+68=\", \"93-45=\"),\n    @(\"30+13=\", \"31+1=\"),\n    @(\"47+28=\", \"28+71=\"),\n    @(\"23+70=\", \"45-28=\"),\n    @(\"52+43=\", \"40+41=\"),\n    @(\"3+93=\", \"34+61=\"),\n    @(\"69-9=\", \"41-25=\"),\n    @(\"35+10=\", \"17+72=\"),\n    @(\"18-16=\", \"48-45=\"),\n    @(\"44-38=\", \"54+14=\"),\n    @(\"11+88=\", \"7+62=\"),\n    @(\"41+42=\", \"50+30=\"),\n    @(\"69-21=\", \"19+65=\"),\n    @(\"53+10=\", \"58+20=\"),\n    @(\"63+1=\", \"53+22=\"),\n    @(\"34-30=\", \"19+50=\"),\n    @(\"17-11=\", \"40+44=\"),\n    @(\"83-59=\", \"61-47=\"),\n    @(\"59-35=\", \"89-53=\"),\n    @(\"92-74=\", \"75-0=\"),\n    @(\"44+22=\", \"72-8=\"),\n    @(\"80-2=\", \"63+29=\"),\n    @(\"17+63=\", \"11+68=\"),\n    @(\"46-37=\", \"91+5=\"),\n    @(\"3+19=\", \"48+31=\"),\n    @(\"19+39=\", \"61-46=\"),\n)\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$total = $paras.Count\n\n$targetIndex = 0\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $paras.Item($i)\n    $r = $p.Range\n    $rawText = $r.Text\n    $cleanText = $rawText -replace \"[`r`a`f]\", \"\"\n    if ($cleanText.Trim() -eq \"\") {\n        continue\n    }\n\n    if ($targetIndex -ge $replacements.Count) {\n        throw \"More non-blank paragraphs than expected replacements\"\n    }\n\n    $pair = $replacements[$targetIndex]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    if ($cleanText -ne $oldText) {\n        throw \"Paragraph at index $i (replacement #$targetIndex) mismatch: expected '$oldText', found '$cleanText'\"\n    }\n\n    $r.Text = $newText\n    $targetIndex = $targetIndex + 1\n}\n\nif ($targetIndex -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) replacements, applied $targetIndex\"\n}\n\nWrite-Output \"Applied $targetIndex replacements\"\n"}
